$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.608.21'
$ws.Range("E2").Value = '  +2.04%  '

# Row 3
$ws.Range("D3").Value = '1.888.28'
$ws.Range("E3").Value = '  +0.26%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.94%  '

# Row 6
$ws.Range("E6").Value = '  +0.21%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4912'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.64%  '

# Row 8
$ws.Range("E8").Value = '  -0.10%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06759'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.65%  '

# Row 10
$ws.Range("D10").Value = '1.885.60'
$ws.Range("E10").Value = '  +0.15%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.16'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.33%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07240'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.65%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '90.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.07%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6758'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.01%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.044'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.15%  '

# Row 16
$ws.Range("D16").Value = '30.567.89'
$ws.Range("E16").Value = '  +2.02%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007945'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.92%  '

# Row 18
$ws.Range("E18").Value = '  +0.22%  '

# Row 19
$ws.Range("E19").Value = '  +2.54%  '

# Row 20
$ws.Range("D20").Value = '2.131.37'
$ws.Range("E20").Value = '  +0.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.16%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.816'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.64%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '183.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +28.61%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.054'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.76%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.328'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.25%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.26%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.899'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.36%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.398'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.16%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.314'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.29%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09037'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.97%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.993'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.58%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05193'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.67%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7499'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.69%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.110'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.68%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01847'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.89%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.657'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.53%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.135'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.02%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9362'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.11%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4410'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.94%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.88%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.731'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.55%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.569'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.07%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1334'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.89%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05849'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.95%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.443'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.47%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.629'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.67%  '

# Row 50
$ws.Range("E50").Value = '  +3.77%  '

# Row 51
$ws.Range("E51").Value = '  +2.39%  '
